$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.390522308914854
$ws.Range("D2").Value = 0.0474846538192395
$ws.Range("E2").Value = 0.06383708548180067
$ws.Range("F2").Value = 7.107209821577129
$ws.Range("G2").Value = 0.002704054131961699
$ws.Range("J2").Value = 0.2222422812520186
$ws.Range("K2").Value = 1.217257839327175
$ws.Range("L2").Value = 0.277192302064627
$ws.Range("M2").Value = 0.3371880716758398
$ws.Range("B3").Value = 1.388129126599779
$ws.Range("D3").Value = 0.04154735179352542
$ws.Range("E3").Value = 0.06316964206907283
$ws.Range("F3").Value = 6.916827027181171
$ws.Range("G3").Value = 0.002709636068185649
$ws.Range("J3").Value = 0.2192234540388043
$ws.Range("K3").Value = 1.200766569114137
$ws.Range("L3").Value = 0.279273351339711
$ws.Range("M3").Value = 0.3376836515447614
$ws.Range("B4").Value = 1.38750667872975
$ws.Range("D4").Value = 0.03789072640053348
$ws.Range("E4").Value = 0.06274793392579081
$ws.Range("F4").Value = 6.800554705658641
$ws.Range("G4").Value = 0.002713241681167302
$ws.Range("J4").Value = 0.2173162972436344
$ws.Range("K4").Value = 1.191894124772347
$ws.Range("L4").Value = 0.2808088131310242
$ws.Range("M4").Value = 0.3382293558468241
$ws.Range("B5").Value = 1.387466182014293
$ws.Range("D5").Value = 0.03639767680942896
$ws.Range("E5").Value = 0.06257306022841647
$ws.Range("F5").Value = 6.753325943201759
$ws.Range("G5").Value = 0.002714755985314774
$ws.Range("J5").Value = 0.2165254835057837
$ws.Range("K5").Value = 1.188592896684753
$ws.Range("L5").Value = 0.281499269593894
$ws.Range("M5").Value = 0.3385124572663258
$ws.Range("B6").Value = 1.387472337026452
$ws.Range("D6").Value = 0.03614957468968782
$ws.Range("E6").Value = 0.06254383896656446
$ws.Range("F6").Value = 6.745492787971102
$ws.Range("G6").Value = 0.002715010155990348
$ws.Range("J6").Value = 0.216393341555019
$ws.Range("K6").Value = 1.188063694810978
$ws.Range("L6").Value = 0.2816178281610391
$ws.Range("M6").Value = 0.3385631340651472
$ws.Range("B7").Value = 1.387505269284048
$ws.Range("D7").Value = 0.0378706027377973
$ws.Range("E7").Value = 0.06274558780551853
$ws.Range("F7").Value = 6.79991714675873
$ws.Range("G7").Value = 0.002713261921445258
$ws.Range("J7").Value = 0.2173056874527823
$ws.Range("K7").Value = 1.191848331432226
$ws.Range("L7").Value = 0.2808178627990685
$ws.Range("M7").Value = 0.3382329279675886
$ws.Range("B8").Value = 1.389521448285791
$ws.Range("D8").Value = 0.0454396808161448
$ws.Range("E8").Value = 0.06360939978419022
$ws.Range("F8").Value = 7.041434180147064
$ws.Range("G8").Value = 0.0027059418760745
$ws.Range("J8").Value = 0.2212124157837749
$ws.Range("K8").Value = 1.211311049954361
$ws.Range("L8").Value = 0.2778563299476531
$ws.Range("M8").Value = 0.3373088496200971
$ws.Range("B9").Value = 1.40019055047344
$ws.Range("D9").Value = 0.06020216905233156
$ws.Range("E9").Value = 0.06521048262967355
$ws.Range("F9").Value = 7.520183765442908
$ws.Range("G9").Value = 0.002692994684308444
$ws.Range("J9").Value = 0.2284555873034328
$ws.Range("K9").Value = 1.259462255140051
$ws.Range("L9").Value = 0.2740963977939543
$ws.Range("M9").Value = 0.3374122118727456
$ws.Range("B10").Value = 1.412120449544005
$ws.Range("D10").Value = 0.07101058087449985
$ws.Range("E10").Value = 0.06633237734482655
$ws.Range("F10").Value = 7.875348223139753
$ws.Range("G10").Value = 0.002684330282109614
$ws.Range("J10").Value = 0.2335326452864983
$ws.Range("K10").Value = 1.300990768085882
$ws.Range("L10").Value = 0.2725871306696348
$ws.Range("M10").Value = 0.3386568066717501
$ws.Range("B11").Value = 1.418435951885584
$ws.Range("D11").Value = 0.07592197860168426
$ws.Range("E11").Value = 0.06683141644823998
$ws.Range("F11").Value = 8.037732498390881
$ws.Range("G11").Value = 0.002680570578876867
$ws.Range("J11").Value = 0.2357914690139822
$ws.Range("K11").Value = 1.321233836183893
$ws.Range("L11").Value = 0.2721737985802832
$ws.Range("M11").Value = 0.3394770835070915
$ws.Range("B12").Value = 1.420955175327322
$ws.Range("D12").Value = 0.07778120601091132
$ws.Range("E12").Value = 0.06701880002348481
$ws.Range("F12").Value = 8.099345554729211
$ws.Range("G12").Value = 0.002679172851587168
$ws.Range("J12").Value = 0.2366397039111909
$ws.Range("K12").Value = 1.329094815519341
$ws.Range("L12").Value = 0.2720566650386758
$ws.Range("M12").Value = 0.3398242570427676
$ws.Range("B13").Value = 1.420406938219003
$ws.Range("D13").Value = 0.07738081272086106
$ws.Range("E13").Value = 0.0669785138389436
$ws.Range("F13").Value = 8.086070605674877
$ws.Range("G13").Value = 0.002679472723414865
$ws.Range("J13").Value = 0.2364573359228608
$ws.Range("K13").Value = 1.327393107369431
$ws.Range("L13").Value = 0.2720801385734859
$ws.Range("M13").Value = 0.3397478610937483
$ws.Range("B14").Value = 1.41864065150213
$ws.Range("D14").Value = 0.07607494968394235
$ws.Range("E14").Value = 0.06684686430260545
$ws.Range("F14").Value = 8.042798980346788
$ws.Range("G14").Value = 0.002680455066893989
$ws.Range("J14").Value = 0.2358613958037985
$ws.Range("K14").Value = 1.321876642270098
$ws.Range("L14").Value = 0.2721633720165073
$ws.Range("M14").Value = 0.339504913033295
$ws.Range("B15").Value = 1.417575375058874
$ws.Range("D15").Value = 0.0752749960347785
$ws.Range("E15").Value = 0.06676601893628131
$ws.Range("F15").Value = 8.016309823214726
$ws.Range("G15").Value = 0.002681060161060593
$ws.Range("J15").Value = 0.2354954413372248
$ws.Range("K15").Value = 1.318523121601316
$ws.Range("L15").Value = 0.2722194870454757
$ws.Range("M15").Value = 0.3393608609719294
$ws.Range("B16").Value = 1.411725588223959
$ws.Range("D16").Value = 0.07068951480383134
$ws.Range("E16").Value = 0.06629953928409016
$ws.Range("F16").Value = 7.864752833648339
$ws.Range("G16").Value = 0.002684579632625052
$ws.Range("J16").Value = 0.2333840186778566
$ws.Range("K16").Value = 1.299695118900985
$ws.Range("L16").Value = 0.2726196491982193
$ws.Range("M16").Value = 0.3386083138882512
$ws.Range("B17").Value = 1.408364440574388
$ws.Range("D17").Value = 0.06787520347539555
$ws.Range("E17").Value = 0.06601049592489883
$ws.Range("F17").Value = 7.771989916119878
$ws.Range("G17").Value = 0.002686785165970294
$ws.Range("J17").Value = 0.2320758473541389
$ws.Range("K17").Value = 1.288491627257315
$ws.Range("L17").Value = 0.2729351853532194
$ws.Range("M17").Value = 0.3382117408486955
$ws.Range("B18").Value = 1.406514834922604
$ws.Range("D18").Value = 0.06625596755517904
$ws.Range("E18").Value = 0.06584317852358712
$ws.Range("F18").Value = 7.718711908092217
$ws.Range("G18").Value = 0.002688070848862469
$ws.Range("J18").Value = 0.2313186344872982
$ws.Range("K18").Value = 1.282174825899688
$ws.Range("L18").Value = 0.2731423862499085
$ws.Range("M18").Value = 0.33800755909866
$ws.Range("B19").Value = 1.40590295795505
$ws.Range("D19").Value = 0.06570762805195329
$ws.Range("E19").Value = 0.0657863433375967
$ws.Range("F19").Value = 7.700685954029183
$ws.Range("G19").Value = 0.002688509103707853
$ws.Range("J19").Value = 0.2310614279555416
$ws.Range("K19").Value = 1.280057873131
$ws.Range("L19").Value = 0.2732169542993361
$ws.Range("M19").Value = 0.3379425340245525
$ws.Range("B20").Value = 1.408713584993251
$ws.Range("D20").Value = 0.06817484419333653
$ws.Range("E20").Value = 0.06604137528930387
$ws.Range("F20").Value = 7.781856712967738
$ws.Range("G20").Value = 0.002686548612629736
$ws.Range("J20").Value = 0.2322155986949426
$ws.Range("K20").Value = 1.28967109115996
$ws.Range("L20").Value = 0.2728989341373307
$ws.Range("M20").Value = 0.3382514813704489
$ws.Range("B21").Value = 1.419155988202704
$ws.Range("D21").Value = 0.07645852861236335
$ws.Range("E21").Value = 0.06688557587578625
$ws.Range("F21").Value = 8.055505571766901
$ws.Range("G21").Value = 0.002680165824776279
$ws.Range("J21").Value = 0.2360366301126327
$ws.Range("K21").Value = 1.32349165128133
$ws.Range("L21").Value = 0.2721378546276725
$ws.Range("M21").Value = 0.3395752806276313
$ws.Range("B22").Value = 1.426724883000333
$ws.Range("D22").Value = 0.08186890354008369
$ws.Range("E22").Value = 0.06742804887621379
$ws.Range("F22").Value = 8.235061917304677
$ws.Range("G22").Value = 0.002676145728599911
$ws.Range("J22").Value = 0.2384923983576428
$ws.Range("K22").Value = 1.346734593555396
$ws.Range("L22").Value = 0.2718700657737543
$ws.Range("M22").Value = 0.3406535239652513
$ws.Range("B23").Value = 1.42261715598184
$ws.Range("D23").Value = 0.07898154936802371
$ws.Range("E23").Value = 0.06713935634526003
$ws.Range("F23").Value = 8.139162923787637
$ws.Range("G23").Value = 0.002678277523253465
$ws.Range("J23").Value = 0.2371854504900526
$ws.Range("K23").Value = 1.334224812452334
$ws.Range("L23").Value = 0.2719919470984991
$ws.Range("M23").Value = 0.3400585448623659
$ws.Range("B24").Value = 1.408555479006679
$ws.Range("D24").Value = 0.06803938053333525
$ws.Range("E24").Value = 0.06602741828999559
$ws.Range("F24").Value = 7.777395770486862
$ws.Range("G24").Value = 0.002686655503104583
$ws.Range("J24").Value = 0.2321524330982001
$ws.Range("K24").Value = 1.289137468641144
$ws.Range("L24").Value = 0.272915242973454
$ws.Range("M24").Value = 0.3382334405080627
$ws.Range("B25").Value = 1.396585647815726
$ws.Range("D25").Value = 0.05621595764796439
$ws.Range("E25").Value = 0.0647870630574916
$ws.Range("F25").Value = 7.390089073966379
$ws.Range("G25").Value = 0.002696347616375497
$ws.Range("J25").Value = 0.2265397992199247
$ws.Range("K25").Value = 1.245360526069675
$ws.Range("L25").Value = 0.2748938575286815
$ws.Range("M25").Value = 0.3371790994860184
